$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 32.3
$ws.Range("F2").Value = 5.383333333333334

# Row 3
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 22
$ws.Range("F3").Value = 3.666666666666667

# Row 4
$ws.Range("E4").Value = 19.1
$ws.Range("F4").Value = 3.183333333333333
$ws.Range("G4").Value = "APROVADO"

# Row 5
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 18
$ws.Range("F5").Value = 3
